# Generate Report for Handback
# The fa9f417a-ea3a-44c9-be15-a1e57a48db16 file has been handed back (in sync
# with en-US) for both the zh-cn and de-de locales. Update the Overview sheet
# and each locale's detail sheet to reflect the new status and handback time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the fa9f417a-ea3a... file.
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: row 3 is the fa9f417a-ea3a... file.
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-17 11:07:58"

# de-de detail sheet: row 3 is the fa9f417a-ea3a... file.
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-17 11:08:05"
